$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.9
$ws.Range("H2").Value = 3.25
$ws.Range("I2").Value = 4.1
$ws.Range("K2").Value = 10
$ws.Range("L2").Value = 1.29
$ws.Range("M2").Value = 3.5
$ws.Range("N2").Value = 1.98
$ws.Range("O2").Value = 1.83
$ws.Range("R2").Value = 1.8
$ws.Range("S2").Value = 1.91
$ws.Range("T2").Value = 7.5
$ws.Range("U2").Value = 9
$ws.Range("W2").Value = 17
$ws.Range("X2").Value = 15
$ws.Range("Y2").Value = 26
$ws.Range("AA2").Value = 6.5
$ws.Range("AB2").Value = 15
$ws.Range("AC2").Value = 51
$ws.Range("AD2").Value = 251
$ws.Range("AE2").Value = 12
$ws.Range("AF2").Value = 21
$ws.Range("AG2").Value = 15
$ws.Range("AH2").Value = 41
$ws.Range("AI2").Value = 34
$ws.Range("AJ2").Value = 41

# Row 5
$ws.Range("G5").Value = 8.25
$ws.Range("H5").Value = 5.3
$ws.Range("I5").Value = 1.25
$ws.Range("T5").Value = 22
$ws.Range("U5").Value = 50
$ws.Range("V5").Value = 22
$ws.Range("W5").Value = 150
$ws.Range("X5").Value = 70
$ws.Range("Y5").Value = 55
$ws.Range("AA5").Value = 9.75
$ws.Range("AB5").Value = 17.5
$ws.Range("AC5").Value = 65
$ws.Range("AD5").Value = 400
$ws.Range("AE5").Value = 7.7
$ws.Range("AF5").Value = 6.1
$ws.Range("AG5").Value = 7.7
$ws.Range("AH5").Value = 6.8
$ws.Range("AJ5").Value = 19.5

# Row 8
$ws.Range("G8").Value = 1.91
$ws.Range("H8").Value = 3.4
$ws.Range("I8").Value = 3.75
$ws.Range("N8").Value = 1.89
$ws.Range("O8").Value = 1.79
$ws.Range("U8").Value = 9.5
$ws.Range("V8").Value = 9
$ws.Range("W8").Value = 17
$ws.Range("AD8").Value = 201
$ws.Range("AE8").Value = 11
$ws.Range("AF8").Value = 17
$ws.Range("AG8").Value = 13
$ws.Range("AI8").Value = 29
$ws.Range("AJ8").Value = 34

# Row 9
$ws.Range("G9").Value = 4.75
$ws.Range("H9").Value = 3.75
$ws.Range("I9").Value = 1.73
$ws.Range("T9").Value = 13
$ws.Range("U9").Value = 23
$ws.Range("AD9").Value = 201
$ws.Range("AF9").Value = 8.5

# Row 10
$ws.Range("G10").Value = 2.4
$ws.Range("I10").Value = 3
$ws.Range("J10").Value = 1.07
$ws.Range("K10").Value = 9
$ws.Range("L10").Value = 1.33
$ws.Range("M10").Value = 3.25
$ws.Range("W10").Value = 23
$ws.Range("X10").Value = 21
$ws.Range("Z10").Value = 9
$ws.Range("AE10").Value = 9
$ws.Range("AH10").Value = 29
$ws.Range("AI10").Value = 23

# Row 11
$ws.Range("K11").Value = 10
$ws.Range("N11").Value = 2
$ws.Range("O11").Value = 1.8
